$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- 1) Append a new trailing-space run to the last existing paragraph
#        ("I made animation that fit the wall run ... icon on the hud.")
#        by replacing the whole paragraph with itself + an extra run.
$lastPara = $d.Paragraphs.Last
$lastRange = $lastPara.Range
$lastXml = "<w:p $wns><w:pPr><w:rPr><w:lang w:val=`"en-US`"/></w:rPr></w:pPr>" + `
  "<w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t>I made animation that fit the wall run to get some extra feedback for the player. The left side wall run animation even adapts to the angle of the player to the wall. I also added some HUD textures and functionality to make the grapple states easier to see for the player. For example, when the grapple is in range, it shows the player an icon on the hud.</w:t></w:r>" + `
  "<w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t xml:space=`"preserve`"> </w:t></w:r>" + `
  "</w:p>"
$lastRange.InsertXML($lastXml)

# Helper: build a simple one-run paragraph (plain "lang en-US" rPr at pPr level)
function New-SimplePara([string]$text, [bool]$preserve = $false) {
    $space = ""
    if ($preserve) { $space = ' xml:space="preserve"' }
    return "<w:p $wns><w:pPr><w:rPr><w:lang w:val=`"en-US`"/></w:rPr></w:pPr>" + `
      "<w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t$space>$text</w:t></w:r></w:p>"
}

# --- 2) 30.11.2024
$r = $d.Content
$r.Collapse(0)
$r.InsertXML((New-SimplePara "30.11.2024"))

# --- 3) Hanging clouds particle effect paragraph (two runs)
$run1 = "I made the hanging clouds particle effect that positions itself in front of the player based on their velocity."
$run2 = " To better indicate speed. I also added sparks that show up when the grapple projectile hits a target to better indicate when and where the player has connected their grapple hook. Also a steam effect when the grapple shooter spawns the projectile just for more visual feedback."
$xml = "<w:p $wns><w:pPr><w:rPr><w:lang w:val=`"en-US`"/></w:rPr></w:pPr>" + `
  "<w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t>$run1</w:t></w:r>" + `
  "<w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t xml:space=`"preserve`">$run2</w:t></w:r>" + `
  "</w:p>"
$r = $d.Content
$r.Collapse(0)
$r.InsertXML($xml)

# --- 4) Shooting mechanic intro paragraph
$r = $d.Content
$r.Collapse(0)
$r.InsertXML((New-SimplePara "I started with working on the shooting mechanic and I have done a lot of groundwork today including:"))

# --- 5) First bullet item: "A new c++ gun" + " class" (two runs, first has lastRenderedPageBreak)
#        Insert as plain paragraph first, then apply default bullet formatting (this is what
#        creates word/numbering.xml + the relationship + content-type override), then fix up
#        the run XML so it matches the lastRenderedPageBreak + split runs shape.
$r = $d.Content
$r.Collapse(0)
$r.InsertXML("<w:p $wns><w:pPr><w:rPr><w:lang w:val=`"en-US`"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t>A new c++ gun class</w:t></w:r></w:p>")
$bulletPara1 = $d.Paragraphs.Last
$bulletPara1.Range.ListFormat.ApplyBulletDefault()
# Now replace its content with the precise two-run version (keeping the numPr/pStyle Word just added)
$bulletPara1 = $d.Paragraphs.Last
$bulletXml1 = "<w:p $wns><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"1`"/></w:numPr><w:rPr><w:lang w:val=`"en-US`"/></w:rPr></w:pPr>" + `
  "<w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:lastRenderedPageBreak/><w:t>A new c++ gun</w:t></w:r>" + `
  "<w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t xml:space=`"preserve`"> class</w:t></w:r>" + `
  "</w:p>"
$bulletPara1.Range.InsertXML($bulletXml1)

# Helper: build a bulleted ListParagraph with a single run
function New-BulletPara([string]$text, [bool]$preserve = $false) {
    $space = ""
    if ($preserve) { $space = ' xml:space="preserve"' }
    return "<w:p $wns><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"1`"/></w:numPr><w:rPr><w:lang w:val=`"en-US`"/></w:rPr></w:pPr>" + `
      "<w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t$space>$text</w:t></w:r></w:p>"
}

# --- 6) Second bullet
$r = $d.Content
$r.Collapse(0)
$r.InsertXML((New-BulletPara "A template animation blueprint that can be reused for other guns"))

# --- 7) Third bullet
$r = $d.Content
$r.Collapse(0)
$r.InsertXML((New-BulletPara "A shoot animation for the revolver and the hands"))

# --- 8) Fourth bullet (trailing space preserved)
$r = $d.Content
$r.Collapse(0)
$r.InsertXML((New-BulletPara "Shooting inputs trigger the shooting animation for the hand and the gun and decrease ammo " $true))

# --- 9) Closing paragraph (two runs)
$run1 = "The goal is to be able to easily add different guns later if there is the time"
$run2 = " For now the shooting cycle need a lot more polishing but the basic functionality and structure stands already."
$xml = "<w:p $wns><w:pPr><w:rPr><w:lang w:val=`"en-US`"/></w:rPr></w:pPr>" + `
  "<w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t>$run1</w:t></w:r>" + `
  "<w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t xml:space=`"preserve`">$run2</w:t></w:r>" + `
  "</w:p>"
$r = $d.Content
$r.Collapse(0)
$r.InsertXML($xml)

"done"
